# Apply the target edits to NORTH_DAKOTA_2017.xlsx (sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Header row: translate Spanish column headers to snake_case machine names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2. Normalize capitalization of "de"/"del"/"los" -> "De"/"Del"/"Los" in
#    a handful of place names scattered through the data rows.
$ws.Range("A15").Value = "Ciudad De México"
$ws.Range("A22").Value = "Estado De México"
$ws.Range("B23").Value = "Ecatepec De Morelos"
$ws.Range("B29").Value = "Coyuca De Catalán"
$ws.Range("B30").Value = "Huitzuco De Los Figueroa"
$ws.Range("B33").Value = "Cuautepec De Hinojosa"
$ws.Range("B35").Value = "Autlán De Navarro"
$ws.Range("B44").Value = "Tepatitlán De Morelos"
$ws.Range("B46").Value = "Zapotlán Del Rey"
$ws.Range("B64").Value = "Oaxaca De Juárez"
$ws.Range("B71").Value = "Tetela De Ocampo"
$ws.Range("B77").Value = "Santa María Del Río"
$ws.Range("B89").Value = "Nochistlán De Mejía"

# 3. Remove the trailing footnote/source blocks. Delete bottom-most block
#    first so the earlier block's row numbers stay valid.
$ws.Range("A476:D480").EntireRow.Delete()
$ws.Range("A94:D98").EntireRow.Delete()
